$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.603.94"
$ws.Range("E2").Value = "  +0.37%  "

# Row 3
$ws.Range("D3").Value = "1.922.30"
$ws.Range("E3").Value = "  -0.94%  "

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").Value = "'247.83"
$ws.Range("E5").Value = "  +1.36%  "

# Row 6
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.07%  "

# Row 7
$ws.Range("D7").Value = "'0.4716"
$ws.Range("E7").Value = "  +0.80%  "

# Row 8
$ws.Range("D8").Value = "'0.2911"
$ws.Range("E8").Value = "  +0.92%  "

# Row 9
$ws.Range("D9").Value = "'0.06772"
$ws.Range("E9").Value = "  +1.02%  "

# Row 10
$ws.Range("D10").Value = "'106.71"
$ws.Range("E10").Value = "  -0.47%  "

# Row 11
$ws.Range("D11").Value = "'18.84"
$ws.Range("E11").Value = "  +0.20%  "

# Row 12
$ws.Range("D12").Value = "1.946.03"
$ws.Range("E12").Value = "  +0.46%  "

# Row 13
$ws.Range("D13").Value = "'0.07735"
$ws.Range("E13").Value = "  +1.05%  "

# Row 14
$ws.Range("D14").Value = "'5.316"
$ws.Range("E14").Value = "  +2.42%  "

# Row 15
$ws.Range("D15").Value = "'0.6743"
$ws.Range("E15").Value = "  +1.50%  "

# Row 16
$ws.Range("D16").Value = "'283.00"
$ws.Range("E16").Value = "  -7.14%  "

# Row 17
$ws.Range("D17").Value = "30.651.68"
$ws.Range("E17").Value = "  +0.42%  "

# Row 18
$ws.Range("D18").Value = "'0.000007609"
$ws.Range("E18").Value = "  -0.14%  "

# Row 19
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  +0.05%  "

# Row 20
$ws.Range("D20").Value = "'12.95"
$ws.Range("E20").Value = "  -1.00%  "

# Row 21
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.185.65"
$ws.Range("E21").Value = "  +0.81%  "

# Row 22
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.552"
$ws.Range("E22").Value = "  +4.29%  "

# Row 23
$ws.Range("D23").Value = "'1.001"

# Row 24
$ws.Range("D24").Value = "'6.435"
$ws.Range("E24").Value = "  +1.40%  "

# Row 25
$ws.Range("D25").Value = "'9.553"
$ws.Range("E25").Value = "  +2.45%  "

# Row 26
$ws.Range("D26").Value = "'165.78"
$ws.Range("E26").Value = "  -1.76%  "

# Row 27
$ws.Range("D27").Value = "'20.38"
$ws.Range("E27").Value = "  -7.85%  "

# Row 28
$ws.Range("D28").Value = "'2.137"
$ws.Range("E28").Value = "  +3.28%  "

# Row 29
$ws.Range("D29").Value = "'0.1066"
$ws.Range("E29").Value = "  -4.30%  "

# Row 30
$ws.Range("D30").Value = "'1.415"
$ws.Range("E30").Value = "  +3.44%  "

# Row 31
$ws.Range("D31").Value = "'4.199"
$ws.Range("E31").Value = "  +1.34%  "

# Row 32
$ws.Range("D32").Value = "'4.085"
$ws.Range("E32").Value = "  +2.55%  "

# Row 33
$ws.Range("D33").Value = "'0.05056"
$ws.Range("E33").Value = "  -0.93%  "

# Row 34
$ws.Range("D34").Value = "'0.7379"
$ws.Range("E34").Value = "  -0.40%  "

# Row 35
$ws.Range("D35").Value = "'1.145"
$ws.Range("E35").Value = "  -1.33%  "

# Row 36
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.737"
$ws.Range("E36").Value = "  -0.52%  "

# Row 37
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "'0.9981"
$ws.Range("E37").Value = "  -0.07%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02035"
$ws.Range("E38").Value = "  -0.14%  "

# Row 39
$ws.Range("D39").Value = "'2.686"
$ws.Range("E39").Value = "  -0.21%  "

# Row 40
$ws.Range("B40").Value = "Quant"
$ws.Range("C40").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D40").Value = "'111.25"
$ws.Range("E40").Value = "  +1.96%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'2.041"
$ws.Range("E41").Value = "  -1.11%  "

# Row 42
$ws.Range("D42").Value = "'0.4518"
$ws.Range("E42").Value = "  +6.17%  "

# Row 43
$ws.Range("D43").Value = "'0.8742"
$ws.Range("E43").Value = "  -1.06%  "

# Row 44
$ws.Range("D44").Value = "'5.912"
$ws.Range("E44").Value = "  +0.65%  "

# Row 45
$ws.Range("D45").Value = "'1.001"
$ws.Range("E45").Value = "  +0.12%  "

# Row 46
$ws.Range("D46").Value = "'67.86"
$ws.Range("E46").Value = "  -2.89%  "

# Row 47
$ws.Range("D47").Value = "'7.328"
$ws.Range("E47").Value = "  +0.66%  "

# Row 48
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").Value = "'49.43"
$ws.Range("E48").Value = "  -9.90%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.367"
$ws.Range("E49").Value = "  +0.75%  "

# Row 50
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.1264"
$ws.Range("E50").Value = "  +2.80%  "

# Row 51
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'35.16"
$ws.Range("E51").Value = "  +0.55%  "
